$d = $word.ActiveDocument

$d.Content.Find.Execute("NIP 198904132015041005", $true, $false, $false, $false, $false,
                         $true, 1, $false, "NIP. 198904132015041005", 2)
